$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Include #0 sheet: update the System URI for CISIS code system ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R01-EnsembleSavoirFaire-CISIS/FHIR/TRE-R01-EnsembleSavoirFaire-CISIS"

# --- Include #1 sheet: update the System URI for CompetenceExclusive code system ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R40-CompetenceExclusive/FHIR/TRE-R40-CompetenceExclusive"
